$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right before the existing
#    "2022-Q1" sheet, so the tab order becomes 总计, 2022-Q3, 2022-Q1.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q3 = $wb.Worksheets.Add($q1)
$q3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the quarterly fund data.
#    Header row + index column reuse the bold/bordered style already
#    used for the same roles on the "总计" sheet (copied over via
#    PasteSpecial so the style entry is shared, not duplicated).
# ------------------------------------------------------------------
$total.Range("B1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Copy() | Out-Null
$q3.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $q3.Cells.Item(1, $col).Value = $headers[$i]
}

$rows = @(
    @("240004", "华宝动力组合混合A", "14.13", "75.08", "2.52", "0.3561", 10),
    @("016257", "华宝动力组合混合C", "2.03", "75.08", "2.52", "0.0512", 10),
    @("009189", "华宝成长策略混合", "1.62", "80.11", "2.32", "0.0376", 10)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $data = $rows[$r]
    $q3.Cells.Item($rowNum, 1).Value = $r
    # Fund code / size / position figures are stored as TEXT in the
    # source workbook (e.g. leading zeros in fund codes must survive),
    # so force them in with a leading apostrophe the way Excel's UI
    # does for "numbers" that should stay text.
    $q3.Cells.Item($rowNum, 2).Value = "'" + $data[0]
    $q3.Cells.Item($rowNum, 3).Value = $data[1]
    $q3.Cells.Item($rowNum, 4).Value = "'" + $data[2]
    $q3.Cells.Item($rowNum, 5).Value = "'" + $data[3]
    $q3.Cells.Item($rowNum, 6).Value = "'" + $data[4]
    $q3.Cells.Item($rowNum, 7).Value = "'" + $data[5]
    $q3.Cells.Item($rowNum, 8).Value = $data[6]
}

# ------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row above the
#    existing Q1 row and fill it in with the Q3 totals.
# ------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# The insert pushed the old row 2 down to row 3, carrying its
# index-column style along; copy that style back onto the new row 2
# before writing values into it.
$total.Range("A3").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.44

# The pre-existing "2022-Q1" row is now row 3; its index counter bumps
# from 0 to 1 to reflect its new position in the list.
$total.Range("A3").Value = 1
